$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: read current F:V content (cols 6..22) for every affected row ---
# (read everything first, then write, so row-rotations/cycles do not clobber source data)
$row6 = @()
for ($c = 6; $c -le 22; $c++) { $row6 += ,($ws.Cells.Item(6, $c).Value2) }
$row7 = @()
for ($c = 6; $c -le 22; $c++) { $row7 += ,($ws.Cells.Item(7, $c).Value2) }
$row8 = @()
for ($c = 6; $c -le 22; $c++) { $row8 += ,($ws.Cells.Item(8, $c).Value2) }
$row9 = @()
for ($c = 6; $c -le 22; $c++) { $row9 += ,($ws.Cells.Item(9, $c).Value2) }
$row10 = @()
for ($c = 6; $c -le 22; $c++) { $row10 += ,($ws.Cells.Item(10, $c).Value2) }
$row16 = @()
for ($c = 6; $c -le 22; $c++) { $row16 += ,($ws.Cells.Item(16, $c).Value2) }
$row17 = @()
for ($c = 6; $c -le 22; $c++) { $row17 += ,($ws.Cells.Item(17, $c).Value2) }
$row26 = @()
for ($c = 6; $c -le 22; $c++) { $row26 += ,($ws.Cells.Item(26, $c).Value2) }
$row27 = @()
for ($c = 6; $c -le 22; $c++) { $row27 += ,($ws.Cells.Item(27, $c).Value2) }
$row34 = @()
for ($c = 6; $c -le 22; $c++) { $row34 += ,($ws.Cells.Item(34, $c).Value2) }
$row35 = @()
for ($c = 6; $c -le 22; $c++) { $row35 += ,($ws.Cells.Item(35, $c).Value2) }
$row36 = @()
for ($c = 6; $c -le 22; $c++) { $row36 += ,($ws.Cells.Item(36, $c).Value2) }
$row38 = @()
for ($c = 6; $c -le 22; $c++) { $row38 += ,($ws.Cells.Item(38, $c).Value2) }
$row39 = @()
for ($c = 6; $c -le 22; $c++) { $row39 += ,($ws.Cells.Item(39, $c).Value2) }
$row45 = @()
for ($c = 6; $c -le 22; $c++) { $row45 += ,($ws.Cells.Item(45, $c).Value2) }
$row46 = @()
for ($c = 6; $c -le 22; $c++) { $row46 += ,($ws.Cells.Item(46, $c).Value2) }
$row64 = @()
for ($c = 6; $c -le 22; $c++) { $row64 += ,($ws.Cells.Item(64, $c).Value2) }
$row65 = @()
for ($c = 6; $c -le 22; $c++) { $row65 += ,($ws.Cells.Item(65, $c).Value2) }
$row67 = @()
for ($c = 6; $c -le 22; $c++) { $row67 += ,($ws.Cells.Item(67, $c).Value2) }
$row68 = @()
for ($c = 6; $c -le 22; $c++) { $row68 += ,($ws.Cells.Item(68, $c).Value2) }
$row69 = @()
for ($c = 6; $c -le 22; $c++) { $row69 += ,($ws.Cells.Item(69, $c).Value2) }

# --- Step 2: write rotated F:V content back ---
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(6, 6 + $i).Value2 = $row7[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(7, 6 + $i).Value2 = $row6[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(8, 6 + $i).Value2 = $row10[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(9, 6 + $i).Value2 = $row8[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(10, 6 + $i).Value2 = $row9[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(16, 6 + $i).Value2 = $row17[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(17, 6 + $i).Value2 = $row16[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(26, 6 + $i).Value2 = $row27[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(27, 6 + $i).Value2 = $row26[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(34, 6 + $i).Value2 = $row35[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(35, 6 + $i).Value2 = $row36[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(36, 6 + $i).Value2 = $row34[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(38, 6 + $i).Value2 = $row39[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(39, 6 + $i).Value2 = $row38[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(45, 6 + $i).Value2 = $row46[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(46, 6 + $i).Value2 = $row45[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(64, 6 + $i).Value2 = $row65[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(65, 6 + $i).Value2 = $row64[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(67, 6 + $i).Value2 = $row68[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(68, 6 + $i).Value2 = $row69[$i] }
for ($i = 0; $i -lt 17; $i++) { $ws.Cells.Item(69, 6 + $i).Value2 = $row67[$i] }

# --- Step 3: append new match rows 73-78 ---
$srcStyleRow = 72

# row 73
$ws.Cells.Item(73, 1).Value2 = 72
$ws.Cells.Item(73, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(73, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(73, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(73, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(73, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(73, 2).Value2 = "denmark"
$ws.Cells.Item(73, 3).Value2 = "3rd-division"
$ws.Cells.Item(73, 4).Value2 = "2023-2024"
$ws.Cells.Item(73, 5).Value2 = 45226.79166666666
$ws.Cells.Item(73, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(73, 6).Value2 = "Holbaek"
$ws.Cells.Item(73, 7).Value2 = 1
$ws.Cells.Item(73, 8).Value2 = "Ishoj"
$ws.Cells.Item(73, 9).Value2 = 4
$ws.Cells.Item(73, 10).Value2 = 2.94
$ws.Cells.Item(73, 11).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(73, 12).Value2 = 3.12
$ws.Cells.Item(73, 13).Value2 = "27/10/2023 18:31"
$ws.Cells.Item(73, 14).Value2 = 3.43
$ws.Cells.Item(73, 15).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(73, 16).Value2 = 3.43
$ws.Cells.Item(73, 17).Value2 = "27/10/2023 18:31"
$ws.Cells.Item(73, 18).Value2 = 2.01
$ws.Cells.Item(73, 19).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(73, 20).Value2 = 2.12
$ws.Cells.Item(73, 21).Value2 = "27/10/2023 18:31"
$ws.Cells.Item(73, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/holbaek-ishoj-if/Yyaccyw9/"

# row 74
$ws.Cells.Item(74, 1).Value2 = 73
$ws.Cells.Item(74, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(74, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(74, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(74, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(74, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(74, 2).Value2 = "denmark"
$ws.Cells.Item(74, 3).Value2 = "3rd-division"
$ws.Cells.Item(74, 4).Value2 = "2023-2024"
$ws.Cells.Item(74, 5).Value2 = 45226.79166666666
$ws.Cells.Item(74, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(74, 6).Value2 = "Young Boys"
$ws.Cells.Item(74, 7).Value2 = 0
$ws.Cells.Item(74, 8).Value2 = "Holstebro"
$ws.Cells.Item(74, 9).Value2 = 2
$ws.Cells.Item(74, 10).Value2 = 1.37
$ws.Cells.Item(74, 11).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(74, 12).Value2 = 1.47
$ws.Cells.Item(74, 13).Value2 = "27/10/2023 18:41"
$ws.Cells.Item(74, 14).Value2 = 4.49
$ws.Cells.Item(74, 15).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(74, 16).Value2 = 4.62
$ws.Cells.Item(74, 17).Value2 = "27/10/2023 18:41"
$ws.Cells.Item(74, 18).Value2 = 5.89
$ws.Cells.Item(74, 19).Value2 = "26/10/2023 07:12"
$ws.Cells.Item(74, 20).Value2 = 5.33
$ws.Cells.Item(74, 21).Value2 = "27/10/2023 18:41"
$ws.Cells.Item(74, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/young-boys-fd-holstebro/G0l2dehF/"

# row 75
$ws.Cells.Item(75, 1).Value2 = 74
$ws.Cells.Item(75, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(75, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(75, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(75, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(75, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(75, 2).Value2 = "denmark"
$ws.Cells.Item(75, 3).Value2 = "3rd-division"
$ws.Cells.Item(75, 4).Value2 = "2023-2024"
$ws.Cells.Item(75, 5).Value2 = 45227.5625
$ws.Cells.Item(75, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(75, 6).Value2 = "VSK Aarhus"
$ws.Cells.Item(75, 7).Value2 = 1
$ws.Cells.Item(75, 8).Value2 = "Lyseng"
$ws.Cells.Item(75, 9).Value2 = 0
$ws.Cells.Item(75, 10).Value2 = 1.4
$ws.Cells.Item(75, 11).Value2 = "27/10/2023 01:43"
$ws.Cells.Item(75, 12).Value2 = 1.64
$ws.Cells.Item(75, 13).Value2 = "28/10/2023 11:40"
$ws.Cells.Item(75, 14).Value2 = 4.4
$ws.Cells.Item(75, 15).Value2 = "27/10/2023 01:43"
$ws.Cells.Item(75, 16).Value2 = 4.14
$ws.Cells.Item(75, 17).Value2 = "28/10/2023 12:10"
$ws.Cells.Item(75, 18).Value2 = 5.29
$ws.Cells.Item(75, 19).Value2 = "27/10/2023 01:43"
$ws.Cells.Item(75, 20).Value2 = 4.26
$ws.Cells.Item(75, 21).Value2 = "28/10/2023 12:10"
$ws.Cells.Item(75, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/vsk-aarhus-lyseng/fim6eF7L/"

# row 76
$ws.Cells.Item(76, 1).Value2 = 75
$ws.Cells.Item(76, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(76, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(76, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(76, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(76, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(76, 2).Value2 = "denmark"
$ws.Cells.Item(76, 3).Value2 = "3rd-division"
$ws.Cells.Item(76, 4).Value2 = "2023-2024"
$ws.Cells.Item(76, 5).Value2 = 45227.58333333334
$ws.Cells.Item(76, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(76, 6).Value2 = "BK Frem"
$ws.Cells.Item(76, 7).Value2 = 0
$ws.Cells.Item(76, 8).Value2 = "Naesby"
$ws.Cells.Item(76, 9).Value2 = 0
$ws.Cells.Item(76, 10).Value2 = 1.78
$ws.Cells.Item(76, 11).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(76, 12).Value2 = 1.63
$ws.Cells.Item(76, 13).Value2 = "28/10/2023 13:52"
$ws.Cells.Item(76, 14).Value2 = 3.56
$ws.Cells.Item(76, 15).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(76, 16).Value2 = 3.92
$ws.Cells.Item(76, 17).Value2 = "28/10/2023 13:52"
$ws.Cells.Item(76, 18).Value2 = 3.49
$ws.Cells.Item(76, 19).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(76, 20).Value2 = 4.59
$ws.Cells.Item(76, 21).Value2 = "28/10/2023 13:52"
$ws.Cells.Item(76, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/bk-frem-naesby/OQiAfZMR/"

# row 77
$ws.Cells.Item(77, 1).Value2 = 76
$ws.Cells.Item(77, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(77, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(77, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(77, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(77, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(77, 2).Value2 = "denmark"
$ws.Cells.Item(77, 3).Value2 = "3rd-division"
$ws.Cells.Item(77, 4).Value2 = "2023-2024"
$ws.Cells.Item(77, 5).Value2 = 45227.58333333334
$ws.Cells.Item(77, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(77, 6).Value2 = "SfB-Oure"
$ws.Cells.Item(77, 7).Value2 = 0
$ws.Cells.Item(77, 8).Value2 = "Vejgaard"
$ws.Cells.Item(77, 9).Value2 = 5
$ws.Cells.Item(77, 10).Value2 = 3.11
$ws.Cells.Item(77, 11).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(77, 12).Value2 = 2.92
$ws.Cells.Item(77, 13).Value2 = "28/10/2023 13:58"
$ws.Cells.Item(77, 14).Value2 = 3.57
$ws.Cells.Item(77, 15).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(77, 16).Value2 = 3.57
$ws.Cells.Item(77, 17).Value2 = "28/10/2023 13:58"
$ws.Cells.Item(77, 18).Value2 = 1.9
$ws.Cells.Item(77, 19).Value2 = "27/10/2023 02:12"
$ws.Cells.Item(77, 20).Value2 = 2.16
$ws.Cells.Item(77, 21).Value2 = "28/10/2023 13:58"
$ws.Cells.Item(77, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/sfb-oure-vejgaard/QZwxjgVr/"

# row 78
$ws.Cells.Item(78, 1).Value2 = 77
$ws.Cells.Item(78, 1).NumberFormat = $ws.Cells.Item($srcStyleRow, 1).NumberFormat
$ws.Cells.Item(78, 1).Font.Bold = $ws.Cells.Item($srcStyleRow, 1).Font.Bold
$ws.Cells.Item(78, 1).HorizontalAlignment = $ws.Cells.Item($srcStyleRow, 1).HorizontalAlignment
$ws.Cells.Item(78, 1).VerticalAlignment = $ws.Cells.Item($srcStyleRow, 1).VerticalAlignment
$ws.Cells.Item(78, 1).Borders.LineStyle = $ws.Cells.Item($srcStyleRow, 1).Borders.LineStyle
$ws.Cells.Item(78, 2).Value2 = "denmark"
$ws.Cells.Item(78, 3).Value2 = "3rd-division"
$ws.Cells.Item(78, 4).Value2 = "2023-2024"
$ws.Cells.Item(78, 5).Value2 = 45228.54166666666
$ws.Cells.Item(78, 5).NumberFormat = $ws.Cells.Item($srcStyleRow, 5).NumberFormat
$ws.Cells.Item(78, 6).Value2 = "Avarta"
$ws.Cells.Item(78, 7).Value2 = 1
$ws.Cells.Item(78, 8).Value2 = "Vanlose"
$ws.Cells.Item(78, 9).Value2 = 1
$ws.Cells.Item(78, 10).Value2 = 2.36
$ws.Cells.Item(78, 11).Value2 = "28/10/2023 02:13"
$ws.Cells.Item(78, 12).Value2 = 2.86
$ws.Cells.Item(78, 13).Value2 = "29/10/2023 12:30"
$ws.Cells.Item(78, 14).Value2 = 3.2
$ws.Cells.Item(78, 15).Value2 = "28/10/2023 02:13"
$ws.Cells.Item(78, 16).Value2 = 3.06
$ws.Cells.Item(78, 17).Value2 = "29/10/2023 12:30"
$ws.Cells.Item(78, 18).Value2 = 2.56
$ws.Cells.Item(78, 19).Value2 = "28/10/2023 02:13"
$ws.Cells.Item(78, 20).Value2 = 2.45
$ws.Cells.Item(78, 21).Value2 = "29/10/2023 12:30"
$ws.Cells.Item(78, 22).Value2 = "https://www.betexplorer.com/football/denmark/3rd-division/avarta-vanlose/8bTtkDpk/"

Write-Output "done"